$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.061.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.468.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.478.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "

$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.908.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.006.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.468.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.409"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0745"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.58%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.794"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "273.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0906"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.44%  "

$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.732.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.16%  "
